# Apply the cryptos.xlsx price/volume refresh described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.889.99"
$ws.Range("E2").Value = "  -2.98%  "
$ws.Range("D3").Value = "1.621.28"
$ws.Range("E3").Value = "  -3.25%  "
$ws.Range("E4").Value = "  +0.36%  "
$ws.Range("E5").Value = "  +0.42%  "
$ws.Range("D6").Value = "'307.28"
$ws.Range("E6").Value = "  -2.09%  "
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("D8").Value = "'0.3811"
$ws.Range("E8").Value = "  -3.23%  "
$ws.Range("D9").Value = "'1.009"
$ws.Range("E9").Value = "  +0.60%  "
$ws.Range("D10").Value = "'49.58"
$ws.Range("E10").Value = "  -5.37%  "
$ws.Range("D11").Value = "'1.349"
$ws.Range("E11").Value = "  -3.00%  "
$ws.Range("D12").Value = "'0.08434"
$ws.Range("E12").Value = "  -2.46%  "
$ws.Range("D13").Value = "'23.63"
$ws.Range("E13").Value = "  -6.08%  "
$ws.Range("D14").Value = "'6.963"
$ws.Range("E14").Value = "  -4.66%  "
$ws.Range("D15").Value = "'0.00001267"
$ws.Range("E15").Value = "  -3.42%  "
$ws.Range("D16").Value = "'7.401"
$ws.Range("E16").Value = "  -4.60%  "
$ws.Range("D17").Value = "1.619.22"
$ws.Range("E17").Value = "  -4.32%  "
$ws.Range("D18").Value = "'92.53"
$ws.Range("E18").Value = "  -1.23%  "
$ws.Range("D19").Value = "'0.06896"
$ws.Range("E19").Value = "  -2.16%  "
$ws.Range("D20").Value = "'19.86"
$ws.Range("E20").Value = "  -3.40%  "
$ws.Range("D21").Value = "'6.824"
$ws.Range("E21").Value = "  -3.31%  "
$ws.Range("D22").Value = "'1.007"
$ws.Range("E22").Value = "  +0.16%  "
$ws.Range("D23").Value = "'13.32"
$ws.Range("E23").Value = "  -4.56%  "
$ws.Range("D24").Value = "23.897.20"
$ws.Range("E24").Value = "  -2.95%  "
$ws.Range("D25").Value = "'2.383"
$ws.Range("E25").Value = "  +1.37%  "
$ws.Range("D26").Value = "'2.853"
$ws.Range("E26").Value = "  +4.94%  "
$ws.Range("D27").Value = "'22.03"
$ws.Range("E27").Value = "  -4.91%  "
$ws.Range("D28").Value = "'157.09"
$ws.Range("E28").Value = "  -3.34%  "
$ws.Range("D29").Value = "'138.47"
$ws.Range("E29").Value = "  -5.78%  "
$ws.Range("D30").Value = "'5.257"
$ws.Range("E30").Value = "  -8.52%  "
$ws.Range("D31").Value = "'7.649"
$ws.Range("E31").Value = "  -3.10%  "
$ws.Range("D32").Value = "'2.460"
$ws.Range("E32").Value = "  -1.97%  "
$ws.Range("D33").Value = "1.805.56"
$ws.Range("E33").Value = "  -3.24%  "
$ws.Range("D34").Value = "'0.07916"
$ws.Range("E34").Value = "  -5.52%  "
$ws.Range("D35").Value = "'0.02867"
$ws.Range("E35").Value = "  -5.55%  "
$ws.Range("D36").Value = "'6.585"
$ws.Range("E36").Value = "  -4.15%  "
$ws.Range("D37").Value = "'0.9512"
$ws.Range("E37").Value = "  -2.49%  "
$ws.Range("D38").Value = "'0.2642"
$ws.Range("E38").Value = "  -6.33%  "
$ws.Range("D39").Value = "'0.09125"
$ws.Range("E39").Value = "  -3.71%  "
$ws.Range("D40").Value = "'10.19"
$ws.Range("E40").Value = "  -2.88%  "
$ws.Range("D41").Value = "'1.412"
$ws.Range("E41").Value = "  -9.28%  "
$ws.Range("B42").Value = "Aptos"
$ws.Range("C42").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D42").Value = "'13.09"
$ws.Range("E42").Value = "  -3.08%  "
$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").Value = "'0.7397"
$ws.Range("E43").Value = "  -6.40%  "
$ws.Range("D44").Value = "'15.85"
$ws.Range("E44").Value = "  -4.13%  "
$ws.Range("D45").Value = "'0.6792"
$ws.Range("E45").Value = "  -4.65%  "
$ws.Range("D46").Value = "'2.426"
$ws.Range("E46").Value = "  -5.14%  "
$ws.Range("D47").Value = "'4.076"
$ws.Range("E47").Value = "  -2.73%  "
$ws.Range("E48").Value = "  +0.13%  "
$ws.Range("D49").Value = "'0.08226"
$ws.Range("E49").Value = "  -4.91%  "
$ws.Range("D50").Value = "'132.36"
$ws.Range("E50").Value = "  -3.50%  "
$ws.Range("D51").Value = "'1.247"
$ws.Range("E51").Value = "  -5.70%  "

# The apostrophe-prefix trick above tags the touched cells with an implicit
# "quote prefix" style; strip it back off so the cells end up with no style
# index, matching the rest of the data rows (only the header row is styled).
$textCoercedRefs = @("D6", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D16", "D18", "D19", "D20", "D21", "D22", "D23", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D49", "D50", "D51")
foreach ($ref in $textCoercedRefs) {
    $ws.Range($ref).Style = "Normal"
}

Write-Output "cryptos list updated"
